$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Row 15 ---
$ws.Range("G15").Value = 1

# --- Row 16 ---
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 69
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = 18.965517241379
$ws.Range("L16").Value = 38
$ws.Range("M16").Value = 130
$ws.Range("N16").Value = -80.563380281690

# --- Row 17 (C17 becomes the text "0" shared string, keeping style 13) ---
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -48.148148148148
$ws.Range("J17").Value = 85
$ws.Range("K17").Value = -29.411764705882
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -15.492957746478
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# --- Row 18 ---
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 27.272727272727
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = -10.101010101010
$ws.Range("L18").Value = 7.228915662650
$ws.Range("M18").Value = 17.105263157894
$ws.Range("N18").Value = -75

# --- Row 19 ---
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 90
$ws.Range("H19").Value = -17.777777777777
$ws.Range("I19").Value = 488
$ws.Range("J19").Value = 454
$ws.Range("K19").Value = 7.488986784140
$ws.Range("L19").Value = 3.389830508474
$ws.Range("M19").Value = 10.657596371882
$ws.Range("N19").Value = -67.004732927653

# --- Row 20 (C20 becomes the text "0" shared string, keeping style 13) ---
$ws.Range("N20").Value = -96.010638297872
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

# --- Row 21 ---
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -15
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 156
$ws.Range("H21").Value = -14.102564102564
$ws.Range("I21").Value = 730
$ws.Range("J21").Value = 717
$ws.Range("K21").Value = 1.813110181311
$ws.Range("L21").Value = 5.950653120464
$ws.Range("M21").Value = 25.429553264604
$ws.Range("N21").Value = -72.463221425877

# --- Row 22 ---
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = -27.272727272727
$ws.Range("I22").Value = 46
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = 24.324324324324
$ws.Range("L22").Value = 39.393939393939
$ws.Range("M22").Value = 58.620689655172

# --- Row 24 ---
$ws.Range("C24").Value = 56
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = -3.448275862068
$ws.Range("F24").Value = 260
$ws.Range("G24").Value = 263
$ws.Range("H24").Value = -1.140684410646
$ws.Range("I24").Value = 1449
$ws.Range("J24").Value = 1663
$ws.Range("K24").Value = -12.868310282621
$ws.Range("L24").Value = -8.638083228247
$ws.Range("M24").Value = 104.950495049505

# --- Row 25 ---
$ws.Range("C25").Value = 56
$ws.Range("D25").Value = 58
$ws.Range("E25").Value = -3.448275862068
$ws.Range("F25").Value = 264
$ws.Range("G25").Value = 260
$ws.Range("H25").Value = 1.538461538461
$ws.Range("I25").Value = 1423
$ws.Range("J25").Value = 1649
$ws.Range("K25").Value = -13.705275924802
$ws.Range("L25").Value = -10.333963453056

# --- Row 26 ---
$ws.Range("C26").Value = 9
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = -14.634146341463
$ws.Range("I26").Value = 195
$ws.Range("J26").Value = 190
$ws.Range("K26").Value = 2.631578947368
$ws.Range("L26").Value = 27.450980392156
$ws.Range("M26").Value = 82.242990654205

# --- Row 27 ---
$ws.Range("G27").Value = 1

# --- Row 28 ---
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("F28").Value = 9
$ws.Range("H28").Value = -10
$ws.Range("I28").Value = 43
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = 2.380952380952
$ws.Range("L28").Value = 19.444444444444

# --- Row 31 restructure ---
# Before: C31=text"0"  D31=1        E31=-100   F31=1        G31=2  H31=-50
# After:  C31=text"0"  D31=text"0"  E31=text"***.*"  F31=text"0"  G31=1  H31=-100
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4104)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
